$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.152.10"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").Value = "1.871.82"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "307.43"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.5053"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("D8").Value = "0.3751"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "0.07157"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "0.8898"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").Value = "20.72"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").Value = "1.877.30"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").Value = "0.07571"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "5.327"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "'89.40"
$ws.Range("D17").Value = "'0.000008506"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D20").Value = "27.207.40"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").Value = "5.082"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "2.118.05"
$ws.Range("D23").Value = "10.61"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "6.495"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "150.82"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "1.846"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "18.02"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "2.093"
$ws.Range("E28").Value = "  -5.82%  "
$ws.Range("D29").Value = "112.93"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "4.769"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").Value = "4.689"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").Value = "0.09006"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "0.05131"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").Value = "3.097"
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("D35").Value = "'0.7440"
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("E36").Value = "  -5.57%  "
$ws.Range("D37").Value = "0.02037"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("D38").Value = "2.544"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "3.045"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "1.077"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").Value = "0.5376"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("D42").Value = "6.608"
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("D43").Value = "114.94"
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").Value = "8.441"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "'0.1480"
$ws.Range("E45").Value = "  -2.65%  "
$ws.Range("D46").Value = "0.4651"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  -5.16%  "
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").Value = "64.71"
$ws.Range("E50").Value = "  -4.11%  "
$ws.Range("D51").Value = "36.55"
$ws.Range("E51").Value = "  -1.68%  "
